$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Occasionally employed.jamais.pro"
$ws.Range("C1").Value = "Regularly employed.jamais.pro"
$ws.Range("D1").Value = "Student.jamais.pro"
$ws.Range("E1").Value = "Unemployed / discouraged.jamais.pro"
$ws.Range("F1").Value = "Receiving social benefits / pensioners / house-makers / disable.jamais.pro"
$ws.Range("G1").Value = "Other.jamais.pro"
$ws.Range("H1").Value = "Not known / missing.jamais.pro"
$ws.Range("I1").Value = "Total.jamais.pro"
